# Applies the "Added new models and updates" revisions to the Financials
# sheet (Adobe.xlsx to-be-modeled workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Financials")

# --- Row 31: Basic Weighted Avg Shares ---
$ws.Range("O31").Value = 460000000
$ws.Range("P31").Value = 467000000

# --- Row 32: Basic EPS, GAAP ---
$ws.Range("O32").Value = 11.13
$ws.Range("P32").Value = 2.52

# --- Row 33: Basic EPS from Cont Ops ---
$ws.Range("O33").Value = 11.13
$ws.Range("P33").Value = 2.52

# --- Row 34: Diluted Weighted Avg Shares ---
$ws.Range("O34").Value = 461250000
$ws.Range("P34").Value = 467000000

# --- Row 35: Diluted EPS, GAAP ---
$ws.Range("O35").Value = 11.1

# --- Row 36: Diluted EPS from Cont Ops ---
$ws.Range("O36").Value = 11.1

# --- Row 44: Sales per Employee - clear B, C, H, I (now blank/not reported) ---
$ws.Range("B44").ClearContents()
$ws.Range("C44").ClearContents()
$ws.Range("H44").ClearContents()
$ws.Range("I44").ClearContents()

# --- Row 122: Shares Outstanding ---
$ws.Range("H122").Value = 494456000

# --- Row 129: Number of Employees - clear B, C, H, I (now blank/not reported) ---
$ws.Range("B129").ClearContents()
$ws.Range("C129").ClearContents()
$ws.Range("H129").ClearContents()
$ws.Range("I129").ClearContents()

# --- Row 136: Deferred Income Taxes ---
$ws.Range("O136").Value = -230000000
$ws.Range("S136").Value = -108000000

# --- Row 138: Other Non-Cash Adj ---
$ws.Range("O138").Value = 65000000
$ws.Range("S138").Value = 16000000

# --- Row 140: (Inc) Dec in Accts Receiv ---
$ws.Range("O140").Value = -127000000
$ws.Range("S140").Value = -168000000

# --- Row 142: (Inc) Dec in Prepaid Assets ---
$ws.Range("O142").Value = -748000000
$ws.Range("S142").Value = -225000000

# --- Row 143: Inc (Dec) in Accts Payable ---
$ws.Range("O143").Value = 865000000
$ws.Range("S143").Value = 193000000

# --- Row 144: Inc (Dec) in Other ---
$ws.Range("O144").Value = 540000000
$ws.Range("S144").Value = 102000000

# --- Row 154: Net Change in LT Investment ---
$ws.Range("O154").Value = 958000000
$ws.Range("S154").Value = 236000000

# --- Row 155: Dec in LT Investment ---
$ws.Range("O155").Value = 1205000000
$ws.Range("S155").Value = 236000000

# --- Row 161: Other Investing Activities - O and S move back to blank ---
$ws.Range("O161").ClearContents()
$ws.Range("S161").ClearContents()

# --- Row 182: Free Cash Flow per Basic Share ---
$ws.Range("O182").Value = 16.58
$ws.Range("P182").Value = 4.78

# --- Row 183: Price/Free Cash Flow ---
$ws.Range("O183").Value = 30.59
$ws.Range("P183").Value = 19.86
$ws.Range("Q183").Value = 18.74
$ws.Range("R183").Value = 23.7
$ws.Range("S183").Value = 30.92
